$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 290.28125
$ws.Range("I33").Value = 109.10714
$ws.Range("K33").Value = 109.10714
$ws.Range("M33").Value = 119.89286
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 649.8333
$ws.Range("I5").Value = 449.875
$ws.Range("J5").Value = 1049.75
$ws.Range("K5").Value = 449.875
$ws.Range("L5").Value = 1049.75
$ws.Range("M5").Value = -337.875
$ws.Range("N5").Value = -1273.75
$ws.Range("H61").Value = 3326.2144
$ws.Range("I61").Value = 2779.6
$ws.Range("K61").Value = 2779.6
$ws.Range("M61").Value = -2567.6
$ws.Range("H74").Value = 2175.8076
$ws.Range("I74").Value = 1919.1875
$ws.Range("J74").Value = 2586.4
$ws.Range("K74").Value = 1919.1875
$ws.Range("L74").Value = 2586.4
$ws.Range("M74").Value = -1045.1875
$ws.Range("N74").Value = -4334.4
$ws.Range("H77").Value = 2175.8076
$ws.Range("I77").Value = 1919.1875
$ws.Range("J77").Value = 2586.4
$ws.Range("K77").Value = 9595.9375
$ws.Range("L77").Value = 12932
$ws.Range("M77").Value = -5227.9375
$ws.Range("N77").Value = -21668
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H136").Value = 3326.2144
$ws.Range("I136").Value = 2779.6
$ws.Range("K136").Value = 8338.799999999999
$ws.Range("M136").Value = -5788.799999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 649.8333
$ws.Range("I4").Value = 449.875
$ws.Range("J4").Value = 1049.75
$ws.Range("K4").Value = 449.875
$ws.Range("L4").Value = 1049.75
$ws.Range("M4").Value = -334.875
$ws.Range("N4").Value = -1279.75
$ws.Range("H75").Value = 9647.4
$ws.Range("I75").Value = 2745.6667
$ws.Range("J75").Value = 20000
$ws.Range("K75").Value = 2745.6667
$ws.Range("L75").Value = 20000
$ws.Range("M75").Value = -1809.6667
$ws.Range("N75").Value = -21872
$ws.Range("H78").Value = 9647.4
$ws.Range("I78").Value = 2745.6667
$ws.Range("J78").Value = 20000
$ws.Range("K78").Value = 8237.000100000001
$ws.Range("L78").Value = 60000
$ws.Range("M78").Value = -3557.000100000001
$ws.Range("N78").Value = -69360
$ws.Range("H134").Value = 2868.9614
$ws.Range("I134").Value = 2686.7827
$ws.Range("K134").Value = 8060.348100000001
$ws.Range("M134").Value = -5525.348100000001
$ws.Range("H135").Value = 63286.855
$ws.Range("J135").Value = 63286.855
$ws.Range("L135").Value = 63286.855
$ws.Range("N135").Value = -73426.85500000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4277120.5
$ws.Range("I31").Value = 1743.12
$ws.Range("J31").Value = 6293808
$ws.Range("K31").Value = 1743.12
$ws.Range("L31").Value = 6293808
$ws.Range("M31").Value = -1448.12
$ws.Range("N31").Value = -6294398
$ws.Range("H34").Value = 4277120.5
$ws.Range("I34").Value = 1743.12
$ws.Range("J34").Value = 6293808
$ws.Range("K34").Value = 1743.12
$ws.Range("L34").Value = 6293808
$ws.Range("M34").Value = -1541.12
$ws.Range("N34").Value = -6294212
$ws.Range("H58").Value = 2814.9443
$ws.Range("I58").Value = 2535.7144
$ws.Range("J58").Value = 2992.6365
$ws.Range("K58").Value = 2535.7144
$ws.Range("L58").Value = 2992.6365
$ws.Range("M58").Value = -2332.7144
$ws.Range("N58").Value = -3398.6365
$ws.Range("H136").Value = 2814.9443
$ws.Range("I136").Value = 2535.7144
$ws.Range("J136").Value = 2992.6365
$ws.Range("K136").Value = 7607.1432
$ws.Range("L136").Value = 8977.9095
$ws.Range("M136").Value = -5057.1432
$ws.Range("N136").Value = -14077.9095
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 1350.5
$ws.Range("J35").Value = 1350.5
$ws.Range("L35").Value = 4051.5
$ws.Range("N35").Value = -4627.5
$ws.Range("H68").Value = 1182.1538
$ws.Range("J68").Value = 1246.5714
$ws.Range("L68").Value = 3739.7142
$ws.Range("N68").Value = -5361.7142
$ws.Range("H71").Value = 1182.1538
$ws.Range("J71").Value = 1246.5714
$ws.Range("L71").Value = 11219.1426
$ws.Range("N71").Value = -19331.1426
$ws.Range("H99").Value = 3910.7144
$ws.Range("I99").Value = 1687.5
$ws.Range("J99").Value = 4800
$ws.Range("K99").Value = 5062.5
$ws.Range("L99").Value = 14400
$ws.Range("M99").Value = -2816.5
$ws.Range("N99").Value = -18892
$ws.Range("H134").Value = 47764652
$ws.Range("I134").Value = 50152532
$ws.Range("J134").Value = 7070
$ws.Range("K134").Value = 150457596
$ws.Range("L134").Value = 21210
$ws.Range("M134").Value = -150452526
$ws.Range("N134").Value = -31350
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2437.5
$ws.Range("I122").Value = 2760
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 8280
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -5830
$ws.Range("N122").Value = -10600
$ws.Range("H132").Value = 33336316
$ws.Range("I132").Value = 66668940
$ws.Range("J132").Value = 3692.8
$ws.Range("K132").Value = 200006820
$ws.Range("L132").Value = 11078.4
$ws.Range("M132").Value = -200004290
$ws.Range("N132").Value = -16138.4
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 1006
$ws.Range("J20").Value = 1006
$ws.Range("L20").Value = 1006
$ws.Range("N20").Value = -1458
$ws.Range("H40").Value = 5257.25
$ws.Range("I40").Value = 4680.9
$ws.Range("J40").Value = 6217.8335
$ws.Range("K40").Value = 4680.9
$ws.Range("L40").Value = 6217.8335
$ws.Range("M40").Value = -4544.9
$ws.Range("N40").Value = -6489.8335
$ws.Range("H46").Value = 2087.25
$ws.Range("I46").Value = 950.5217
$ws.Range("J46").Value = 4098.385
$ws.Range("K46").Value = 950.5217
$ws.Range("L46").Value = 4098.385
$ws.Range("M46").Value = -762.5217
$ws.Range("N46").Value = -4474.385
$ws.Range("H68").Value = 3690.2188
$ws.Range("I68").Value = 3623.48
$ws.Range("J68").Value = 3928.5715
$ws.Range("K68").Value = 3623.48
$ws.Range("L68").Value = 3928.5715
$ws.Range("M68").Value = -2874.48
$ws.Range("N68").Value = -5426.5715
$ws.Range("H71").Value = 3690.2188
$ws.Range("I71").Value = 3623.48
$ws.Range("J71").Value = 3928.5715
$ws.Range("K71").Value = 18117.4
$ws.Range("L71").Value = 19642.8575
$ws.Range("M71").Value = -14373.4
$ws.Range("N71").Value = -27130.8575
$ws.Range("H132").Value = 3278.2917
$ws.Range("I132").Value = 2053.1538
$ws.Range("J132").Value = 4726.1816
$ws.Range("K132").Value = 6159.4614
$ws.Range("L132").Value = 14178.5448
$ws.Range("M132").Value = -3629.4614
$ws.Range("N132").Value = -19238.5448
$ws.Range("H133").Value = 45000
$ws.Range("J133").Value = 45000
$ws.Range("L133").Value = 45000
$ws.Range("N133").Value = -50060
$ws.Range("H135").Value = 58463.332
$ws.Range("J135").Value = 58463.332
$ws.Range("L135").Value = 58463.332
$ws.Range("N135").Value = -68603.33199999999
$ws.Range("H136").Value = 2501.8572
$ws.Range("I136").Value = 1893.091
$ws.Range("J136").Value = 3171.5
$ws.Range("K136").Value = 5679.272999999999
$ws.Range("L136").Value = 9514.5
$ws.Range("M136").Value = -3129.272999999999
$ws.Range("N136").Value = -14614.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2857.2856
$ws.Range("I62").Value = 2334
$ws.Range("K62").Value = 2334
$ws.Range("M62").Value = -1710
$ws.Range("H65").Value = 2857.2856
$ws.Range("I65").Value = 2334
$ws.Range("K65").Value = 11670
$ws.Range("M65").Value = -8550
$ws.Range("H123").Value = 35250
$ws.Range("J123").Value = 35250
$ws.Range("L123").Value = 35250
$ws.Range("N123").Value = -45050
$ws.Range("H132").Value = 967884.4
$ws.Range("I132").Value = 1450463.4
$ws.Range("J132").Value = 2726.4
$ws.Range("K132").Value = 4351390.199999999
$ws.Range("L132").Value = 8179.200000000001
$ws.Range("M132").Value = -4348860.199999999
$ws.Range("N132").Value = -13239.2
$ws.Range("H136").Value = 730304.9399999999
$ws.Range("I136").Value = 1229067.9
$ws.Range("J136").Value = 1343.6923
$ws.Range("K136").Value = 3687203.7
$ws.Range("L136").Value = 4031.0769
$ws.Range("M136").Value = -3684653.7
$ws.Range("N136").Value = -9131.0769
